$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Edit 1 - Slide 1, shape "Прямоугольник 9": merge the two runs
# "Выполнить " + "домашнее задание прошлого урока :-Р" into a single run
# (same concatenated text, so re-assigning the full range to itself lets
# the engine coalesce the runs while keeping the original formatting,
# including the Wingdings <a:sym> that was on the second run).
# ---------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$sh1 = $s1.Shapes.Item(4)
$tr1 = $sh1.TextFrame.TextRange
$para1 = $tr1.Paragraphs(3)
$full1 = $para1.Characters(1, $para1.Length - 1)
$full1.Text = $full1.Text

# ---------------------------------------------------------------------
# Edit 2 - Slide 11, shape "Прямоугольник 6"
# ---------------------------------------------------------------------
$s11 = $p.Slides.Item(11)
$sh11 = $s11.Shapes.Item(3)
$tr11 = $sh11.TextFrame.TextRange

# 2a) "2) Посчитать среднюю задержку рейса" -> "...задержку прибытия рейса"
$para5 = $tr11.Paragraphs(5)
$para5.Text = "2) Посчитать среднюю задержку прибытия рейса"

# 2b) "4) " + "Построить " + "график " + "зависимости ... день" -> one run
$para7 = $tr11.Paragraphs(7)
$full7 = $para7.Characters(1, $para7.Length - 1)
$full7.Text = $full7.Text

# 2c) append "/revoscaler" onto the end of the help URL (same run)
$para9 = $tr11.Paragraphs(9)
$urlRun = $para9.Characters(9, $para9.Length - 8)
$urlRun.Text = $urlRun.Text + "/revoscaler"

# ---------------------------------------------------------------------
# Edit 3 - Slide 6, shape "Прямоугольник 9"
# ---------------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$sh6 = $s6.Shapes.Item(4)

# 3a) grow the shape's height (300524,1030293 / 7752907 stay unchanged)
$sh6.Height = 319.89379

$tr6 = $sh6.TextFrame.TextRange
$para4 = $tr6.Paragraphs(4)

# 3b) split off the trailing period, merge the remaining text into one
# run, then append ". " + the new sentence as their own runs.
$periodRange = $para4.Characters(88, 1)
$periodRange.Text = ""
$full4 = $para4.Characters(1, $para4.Length - 1)
$full4.Text = $full4.Text
$para4.InsertAfter(".") | Out-Null
$para4.InsertAfter(" ") | Out-Null
$para4.InsertAfter("Учитывать 20 наиболее продаваемых позиций.") | Out-Null

# 3c) add a brand-new bulleted paragraph after it
$newPara = $para4.InsertAfter("`rПостроить scatterplot и линейную модель зависимости кредитного лимита заказчика от количества его заказов.")

$para5b = $tr6.Paragraphs(5)
$subScatter = $para5b.Characters(11, 12)
$subScatter.Text = $subScatter.Text
